# [Jallal] more code clean up
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 9-12 with new address/route data, replacing the old
# repeated rows that used to sit there.
$ws.Range("A9").Value = "2175 E JOY RD, ANN ARBOR, MI 48105-9230"
$ws.Range("B9").Value = "48105-R007"

$ws.Range("A10").Value = "2876 BUTTERNUT ST, ANN ARBOR, MI 48108-1851"
$ws.Range("B10").Value = "48108-C024"

$ws.Range("A11").Value = "1793 ADDINGTON LN, ANN ARBOR, MI 48108-8956"
$ws.Range("B11").Value = "48108-R015"

$ws.Range("A12").Value = "1698 POND SHORE DR, ANN ARBOR, MI 48108-9566"
$ws.Range("B12").Value = "48108-R005"

# Remove the trailing duplicate rows 13-15 that are no longer needed.
$ws.Range("A13:B15").ClearContents()
